$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "60.609.04"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -3.25%  "

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.344.76"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -2.96%  "

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "567.01"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.18%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "146.61"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.20%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.32%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.90"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.32%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -1.49%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.415"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.19%  "

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.912.48"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -3.04%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +1.01%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "27.71"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.14%  "

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.343.32"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.76%  "

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000168"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -2.04%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "60.596.19"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.34%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -1.01%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.94%  "

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.90"
$cell.Style = "Normal"

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "376.86"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.56%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.89%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "74.72"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.85%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.02%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.489.86"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.57%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -5.86%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -4.85%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.73%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.35"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -4.02%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -1.56%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.00%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -4.11%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.91"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.35%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -3.86%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -1.51%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -4.99%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Aptos"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.82"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.93%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Monero"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "167.31"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.27%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.00"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -12.77%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.379.47"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.90%  "

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0747"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -3.56%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -3.66%  "

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.28"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.94%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -3.23%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -4.99%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.454.33"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -4.63%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -3.39%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.03%  "

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.35"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.57%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -2.14%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.12%  "
